# Scheduled market-data refresh: updates computed price/profit columns
# (currentAveragePrice[/NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) on each
# sheet with freshly-fetched values for the affected Leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 5243.643
$ws.Range("I132").Value = 4905.3477
$ws.Range("J132").Value = 6799.8
$ws.Range("K132").Value = 14716.0431
$ws.Range("L132").Value = 20399.4
$ws.Range("M132").Value = -12186.0431
$ws.Range("N132").Value = -25459.4

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 8233.333000000001
$ws.Range("I135").Value = 4750
$ws.Range("J135").Value = 9975
$ws.Range("K135").Value = 42750
$ws.Range("L135").Value = 89775
$ws.Range("M135").Value = -40215
$ws.Range("N135").Value = -94845

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2391.7234
$ws.Range("I138").Value = 2953.3076
$ws.Range("J138").Value = 2301.5925
$ws.Range("K138").Value = 8859.9228
$ws.Range("L138").Value = 6904.7775
$ws.Range("M138").Value = -3719.9228
$ws.Range("N138").Value = -17184.7775

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 6285.5293
$ws.Range("I141").Value = 2672.4167
$ws.Range("J141").Value = 14957
$ws.Range("K141").Value = 8017.250100000001
$ws.Range("L141").Value = 44871
$ws.Range("M141").Value = -2837.250100000001
$ws.Range("N141").Value = -55231

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 363982.34
$ws.Range("I32").Value = 407131.25
$ws.Range("J32").Value = 18791.1
$ws.Range("K32").Value = 407131.25
$ws.Range("L32").Value = 18791.1
$ws.Range("M32").Value = -406844.25
$ws.Range("N32").Value = -19365.1

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 125002936
$ws.Range("I86").Value = 125002936
$ws.Range("K86").Value = 125002936
$ws.Range("M86").Value = -125001813

# Row 87: Winter Weather Conditions / Adamantite Dolabra
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 125002936
$ws.Range("I89").Value = 125002936
$ws.Range("K89").Value = 625014680
$ws.Range("M89").Value = -625009064

# Row 90: The Nightsoil Is Dark and Full of Terrors (L) / Adamantite Dolabra
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 32204.656
$ws.Range("I107").Value = 42543.207
$ws.Range("J107").Value = 1189
$ws.Range("K107").Value = 42543.207
$ws.Range("L107").Value = 1189
$ws.Range("M107").Value = -40623.207
$ws.Range("N107").Value = -5029

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 81.63636
$ws.Range("I7").Value = 60.857143
$ws.Range("J7").Value = 118
$ws.Range("K7").Value = 60.857143
$ws.Range("L7").Value = 118
$ws.Range("M7").Value = 52.142857
$ws.Range("N7").Value = -344

# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 966.6667
$ws.Range("I16").Value = 1100
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 1100
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -813
$ws.Range("N16").Value = -1274

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 1729.5
$ws.Range("I58").Value = 1620
$ws.Range("J58").Value = 1802.5
$ws.Range("K58").Value = 1620
$ws.Range("L58").Value = 1802.5
$ws.Range("M58").Value = -1417
$ws.Range("N58").Value = -2208.5

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 1808.2963
$ws.Range("I99").Value = 1417.7142
$ws.Range("J99").Value = 1945
$ws.Range("K99").Value = 1417.7142
$ws.Range("L99").Value = 1945
$ws.Range("M99").Value = 80.28580000000011
$ws.Range("N99").Value = -4941

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 966.6667
$ws.Range("I113").Value = 1100
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 1100
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 1070
$ws.Range("N113").Value = -5040

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 1808.2963
$ws.Range("I126").Value = 1417.7142
$ws.Range("J126").Value = 1945
$ws.Range("K126").Value = 4253.142599999999
$ws.Range("L126").Value = 5835
$ws.Range("M126").Value = -1783.142599999999
$ws.Range("N126").Value = -10775

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 1729.5
$ws.Range("I136").Value = 1620
$ws.Range("J136").Value = 1802.5
$ws.Range("K136").Value = 4860
$ws.Range("L136").Value = 5407.5
$ws.Range("M136").Value = -2310
$ws.Range("N136").Value = -10507.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 382.73334
$ws.Range("I5").Value = 382.73334
$ws.Range("K5").Value = 1148.20002
$ws.Range("M5").Value = -1036.20002

# Row 11: Putting the Squeeze On / Orange Juice
$ws.Range("H11").Value = 199.66667
$ws.Range("I11").Value = 199.66667
$ws.Range("K11").Value = 599.00001
$ws.Range("M11").Value = -459.00001

# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 1141.9474
$ws.Range("I113").Value = 572
$ws.Range("K113").Value = 1716
$ws.Range("M113").Value = 454

# Row 120: A Happy End / Paella
$ws.Range("H120").Value = 3000
$ws.Range("I120").Value = 3000
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 9000
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -4162
$ws.Range("N120").ClearContents()

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 382.73334
$ws.Range("I135").Value = 382.73334
$ws.Range("K135").Value = 3444.60006
$ws.Range("M135").Value = -909.6000599999998

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 94.8
$ws.Range("I2").Value = 43.333332
$ws.Range("K2").Value = 43.333332
$ws.Range("M2").Value = 69.666668

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 5722.7144
$ws.Range("I70").Value = 5752.478
$ws.Range("J70").Value = 5665.6665
$ws.Range("K70").Value = 5752.478
$ws.Range("L70").Value = 5665.6665
$ws.Range("M70").Value = -5482.478
$ws.Range("N70").Value = -6205.6665

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 5722.7144
$ws.Range("I73").Value = 5752.478
$ws.Range("J73").Value = 5665.6665
$ws.Range("K73").Value = 5752.478
$ws.Range("L73").Value = 5665.6665
$ws.Range("M73").Value = -4816.478
$ws.Range("N73").Value = -7537.6665

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1206.9
$ws.Range("I102").Value = 1152.7142
$ws.Range("J102").Value = 1333.3334
$ws.Range("K102").Value = 1152.7142
$ws.Range("L102").Value = 1333.3334
$ws.Range("M102").Value = 469.2858000000001
$ws.Range("N102").Value = -4577.3334

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 224
$ws.Range("I107").Value = 183.66667
$ws.Range("K107").Value = 183.66667
$ws.Range("M107").Value = 1736.33333

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 11138.632
$ws.Range("I22").Value = 1350
$ws.Range("J22").Value = 12290.235
$ws.Range("K22").Value = 1350
$ws.Range("L22").Value = 12290.235
$ws.Range("M22").Value = -1055
$ws.Range("N22").Value = -12880.235

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 11138.632
$ws.Range("I27").Value = 1350
$ws.Range("J27").Value = 12290.235
$ws.Range("K27").Value = 1350
$ws.Range("L27").Value = 12290.235
$ws.Range("M27").Value = -1243
$ws.Range("N27").Value = -12504.235

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 13891846
$ws.Range("I136").Value = 4950
$ws.Range("K136").Value = 14850
$ws.Range("M136").Value = -12300

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 3800.5908
$ws.Range("I136").Value = 3970.2307
$ws.Range("J136").Value = 3555.5557
$ws.Range("K136").Value = 11910.6921
$ws.Range("L136").Value = 10666.6671
$ws.Range("M136").Value = -9360.6921
$ws.Range("N136").Value = -15766.6671

Write-Host "Anima_Profits market data refreshed."